$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 6 (c83b5c62 file) moves from "Ready for handoff"
#     to "Handed back: in sync with en-US" for both zh-cn and de-de columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E6").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F6").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 6 (c83b5c62 file) is handed back.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C6").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("I6").Value = "c83b5c62-0937-4a4b-8082-f6315878c0be.md"
$wsZhCn.Range("J6").Value = "c83b5c62-0937-4a4b-8082-f6315878c0be.210d8acff2aa06cfadd974fe4be4b3451fec0aa9.zh-cn.xlf"
$wsZhCn.Range("K6").Value = "2016-08-31 12:33:56"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/210d8acff2aa06cfadd974fe4be4b3451fec0aa9/e2e/c83b5c62-0937-4a4b-8082-f6315878c0be.md", [Type]::Missing, [Type]::Missing, "c83b5c62-0937-4a4b-8082-f6315878c0be.md")

# --- de-de sheet: row 6 (c83b5c62 file) is handed back.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C6").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("I6").Value = "c83b5c62-0937-4a4b-8082-f6315878c0be.md"
$wsDeDe.Range("J6").Value = "c83b5c62-0937-4a4b-8082-f6315878c0be.210d8acff2aa06cfadd974fe4be4b3451fec0aa9.de-de.xlf"
$wsDeDe.Range("K6").Value = "2016-08-31 12:34:23"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/210d8acff2aa06cfadd974fe4be4b3451fec0aa9/e2e/c83b5c62-0937-4a4b-8082-f6315878c0be.md", [Type]::Missing, [Type]::Missing, "c83b5c62-0937-4a4b-8082-f6315878c0be.md")
